$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing data cells (rows 23 & 25) ---
# D23 is a new value cell (style matches the existing "start/end" time columns, s=2)
$ws.Range("C21").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 0.85416666666666663

# New row 25: date (s=1), start (s=2), end (s=2)
$ws.Range("B21").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = 42926

$ws.Range("C21").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 0.83333333333333337

$ws.Range("C21").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 0.95833333333333337

# --- Header for the new "sum" column ---
$ws.Range("G1").Value = "sum"

# --- Column E: per-row "end minus start" duration ---
$ws.Range("E2").Formula = "=D2-C2"
$ws.Range("E3").Formula = "=D3-C3"
$ws.Range("E4:E25").Formula = "=D4-C4"

# Make the whole E column use the same time format as columns C/D (s=2)
$ws.Range("C21").Copy()
$ws.Range("E2:E25").PasteSpecial(-4122)

# --- G2: grand total of column E, formatted as elapsed time ---
$ws.Range("G2").Formula = "=SUM(E:E)"
$ws.Range("G2").NumberFormat = "[h]:mm:ss"

# --- Selection, matching the post-edit workbook state ---
[void]$ws.Range("F25").Select()
